# Start Yr Shares of Capacity Costs that are Soft Costs.xlsx
#
# This script:
#  1. On the "About" sheet, inserts a new row (26) so the note that used to
#     read as a single line ("Data in this variable are only used for
#     onshore wind, offshore wind, and solar PV.") is split across two
#     lines/rows: "...solar PV" (no trailing period) followed by a new
#     second line "(as well as distributed solar)."
#  2. Updates the second title cell (A2) on the "About" sheet from the old
#     "SYSoCCtaSC Share of Distributed Solar and Retrofitting Costs that is
#     Labor" to the new "SYSoCCtaSC Start Year Share of Distributed Solar
#     Costs that are Soft Costs".
#
# (The other worksheets -- SYSoCCtaSC-electricity and SYSoCCtaSC-buildings --
# only reference the shared strings whose table position shifts as a result
# of the above; their own text/values/formulas are unchanged, so Excel will
# automatically keep pointing at the correct (renumbered) shared string
# entries -- no direct edits are required there.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Insert a new blank row at 26, shifting rows 26-37 down to 27-38, so that
# the note can be split onto two lines (row 25 and new row 26).
$ws.Rows("26:26").Insert()

# Row 25 keeps the first part of the note, minus the trailing period (the
# sentence now continues onto row 26).
$ws.Range("A25").Value = "Data in this variable are only used for onshore wind, offshore wind, and solar PV"

# New row 26 carries the rest of the sentence.
$ws.Range("A26").Value = "(as well as distributed solar)."

# Update the second page title (A2) to the new variable name/description.
$ws.Range("A2").Value = "SYSoCCtaSC Start Year Share of Distributed Solar Costs that are Soft Costs"
